$wb = $excel.ActiveWorkbook

# "Metadata" sheet: row 4 is the "Name" property -> set its Value to "CompetenceVs"
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B4").Value = "CompetenceVs"

# Update the "Date" property value (row 8, column B) to the new generation timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
